# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Rows 48/49 also swap coin identity (Aave <-> FraxShare) to match the
# upstream ranking source for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.847.07"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").Value = "2.673.50"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.77"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.49"
$ws.Range("E6").Value = "  -6.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.598"
$ws.Range("E7").Value = "  -3.01%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("E9").Value = "  -3.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.38"
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.13"
$ws.Range("E12").Value = "  -4.48%  "
$ws.Range("D13").Value = "3.078.00"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "2.675.51"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.930"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.14"
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").Value = "45.909.33"
$ws.Range("E18").Value = "  -4.28%  "
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.85"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.95"
$ws.Range("E21").Value = "  -3.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.23"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "283.21"
$ws.Range("E23").Value = "  +4.36%  "
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "30.95"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.55"
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.44"
$ws.Range("E30").Value = "  -5.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.17"
$ws.Range("E31").Value = "  -6.20%  "
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.61"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0842"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.68"
$ws.Range("E39").Value = "  +8.68%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.24"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.61"
$ws.Range("E42").Value = "  -4.77%  "
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("E44").Value = "  -7.65%  "
$ws.Range("D45").Value = "2.160.01"
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.56"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.35"
$ws.Range("E48").Value = "  -10.15%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "111.86"
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("D50").Value = "2.927.94"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("E51").Value = "  -2.30%  "
